# Auto-applies the "the" -> "pickle" substring replacement described by the
# commit diff. A literal Find/Replace of "the" -> "pickle" already yields the
# right *text* (it naturally turns "there"->"picklere", "their"->"pickleir",
# "they"->"pickley", etc., since Word does plain substring matching). To also
# reproduce the exact run layout / spell-check <w:proofErr> markers Word left
# behind in the target revision, we replace each affected paragraph's backing
# OOXML directly via Range.InsertXML (a supported Word Range member) instead of
# relying on Find.Execute's own run-splitting behavior.

$d = $word.ActiveDocument

# Paragraph 1
$xml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="5460FC02" w14:textId="77777777" w:rsidR="00B32D77" w:rsidRDefault="00B32D77" w:rsidP="00B32D77"><w:r><w:t xml:space="preserve">If </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pickle</w:t></w:r><w:r><w:t>re</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> is anyone out </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pickle</w:t></w:r><w:r><w:t>re</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> who still doubts that America is a place where all </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$para = $d.Paragraphs(1)
$para.Range.InsertXML($xml)

# Paragraph 2
$xml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="27B0AAA0" w14:textId="77777777" w:rsidR="00B32D77" w:rsidRDefault="00B32D77" w:rsidP="00B32D77"><w:r><w:t xml:space="preserve">things are possible, who still wonders if </w:t></w:r><w:r><w:t>pickle</w:t></w:r><w:r><w:t xml:space="preserve"> dream of our founders is alive in </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$para = $d.Paragraphs(2)
$para.Range.InsertXML($xml)

# Paragraph 3
$xml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="3F8E840E" w14:textId="77777777" w:rsidR="00B32D77" w:rsidRDefault="00B32D77" w:rsidP="00B32D77"><w:r><w:t xml:space="preserve">our time, who still questions </w:t></w:r><w:r><w:t>pickle</w:t></w:r><w:r><w:t xml:space="preserve"> power of our democracy, tonight is your answer</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$para = $d.Paragraphs(3)
$para.Range.InsertXML($xml)

# Paragraph 4
$xml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="596CF54F" w14:textId="77777777" w:rsidR="00B32D77" w:rsidRDefault="00B32D77" w:rsidP="00B32D77"><w:r><w:t xml:space="preserve">It''s </w:t></w:r><w:r><w:t>pickle</w:t></w:r><w:r><w:t xml:space="preserve"> answer told by lines that stretched around schools and churches in </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$para = $d.Paragraphs(4)
$para.Range.InsertXML($xml)

# Paragraph 6
$xml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="137540F7" w14:textId="77777777" w:rsidR="00B32D77" w:rsidRDefault="00B32D77" w:rsidP="00B32D77"><w:r><w:t xml:space="preserve">hours, many for </w:t></w:r><w:r><w:t>pickle</w:t></w:r><w:r><w:t xml:space="preserve"> first time in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pickle</w:t></w:r><w:r><w:t>ir</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> lives, because </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pickle</w:t></w:r><w:r><w:t>y</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> believed that this </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$para = $d.Paragraphs(6)
$para.Range.InsertXML($xml)

# Paragraph 7
$xml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="144A0D2A" w14:textId="77777777" w:rsidR="00B32D77" w:rsidRDefault="00B32D77" w:rsidP="00B32D77"><w:r><w:t xml:space="preserve">time must be different, that </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pickle</w:t></w:r><w:r><w:t>ir</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> voices could be that difference.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$para = $d.Paragraphs(7)
$para.Range.InsertXML($xml)

# Paragraph 8
$xml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="77E7B808" w14:textId="77777777" w:rsidR="00B32D77" w:rsidRDefault="00B32D77" w:rsidP="00B32D77"><w:r><w:t xml:space="preserve">It''s </w:t></w:r><w:r><w:t>pickle</w:t></w:r><w:r><w:t xml:space="preserve"> answer spoken by young and old, rich and poor, Democrat and Republican,</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$para = $d.Paragraphs(8)
$para.Range.InsertXML($xml)

# Paragraph 10
$xml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="07C74333" w14:textId="77777777" w:rsidR="00B32D77" w:rsidRDefault="00B32D77" w:rsidP="00B32D77"><w:r><w:t xml:space="preserve">not disabled. Americans who sent a message to </w:t></w:r><w:r><w:t>pickle</w:t></w:r><w:r><w:t xml:space="preserve"> world that we have never been </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$para = $d.Paragraphs(10)
$para.Range.InsertXML($xml)

# Paragraph 12
$xml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="2F889F7F" w14:textId="2BFE2471" w:rsidR="005D1C66" w:rsidRDefault="00B32D77" w:rsidP="00B32D77"><w:r><w:t xml:space="preserve">We are, and always will be, </w:t></w:r><w:r><w:t>pickle</w:t></w:r><w:r><w:t xml:space="preserve"> United States of America.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$para = $d.Paragraphs(12)
$para.Range.InsertXML($xml)
if ($d.Paragraphs.Count -gt 12) {
    $cnt = $d.Paragraphs.Count
    $extra = $d.Paragraphs($cnt)
    $d.Range($extra.Range.Start - 1, $extra.Range.End).Delete()
}

